$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$I = @(2,8,7,9,7,8,8,7,6,7,11,6,7,6,7,6,7,6,8,8,7,5,7,6,7,9,9,8,9,9,6,6,6,7,7,4,7,9,6,5,9,5,8,6,6,7,8,6,9,9,5,9,7,6,8,7,9,4,2,6,8,7,6,4,4)
$J = @(4,8,7,9,7,8,8,8,6,8,11,7,8,6,8,6,7,7,8,8,7,6,8,7,7,9,9,8,9,9,7,6,7,7,7,5,7,9,6,6,10,6,8,7,6,7,8,7,9,9,6,9,7,6,8,7,9,4,2,7,8,7,6,4,4)

for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 9).Value = $I[$r-2]
    $ws.Cells.Item($r, 10).Value = $J[$r-2]
}

$wb.Save()
